# Update presentation with citation
#
# Adds three small "[1]" citation text boxes to the "Confusion Matrix"
# slide (SlideID 266): a full source-link caption at the bottom of the
# slide plus two short "[1]" footnote markers next to the figures that
# the citation refers to.
#
# Shape Left/Top/Width/Height on the PowerPoint COM object model are in
# points (1 pt = 12700 EMU) and are stored as single-precision floats,
# so the literals below are chosen so that they round-trip through that
# float32 conversion to the exact target EMU values.

$p = $ppt.ActivePresentation

# Locate the slide with SlideID 266 (the "Confusion Matrix" slide) rather
# than assuming a fixed index.
$s = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    if ($candidate.SlideID -eq 266) {
        $s = $candidate
        break
    }
}
if ($s -eq $null) {
    $s = $p.Slides.Item(4)
}

# Use the existing "TextBox 11" shape (the unformatted caption above the
# confusion-matrix picture) as a template so the new boxes inherit the
# same body/text formatting (wrap="square" rtlCol="0", spAutoFit,
# noFill, lstStyle, en-CA language) that PowerPoint normally stamps on
# hand-authored text boxes.
$template = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "TextBox 11") {
        $template = $candidate
        break
    }
}

function Add-CitationBox($Name, $Text, $Left, $Top, $Width, $Height) {
    $range = $template.Duplicate()
    $shape = $range.Item(1)
    $shape.Name = $Name
    $shape.TextFrame.TextRange.Text = $Text
    $shape.TextFrame.TextRange.Font.Size = 12
    $shape.Left = $Left
    $shape.Top = $Top
    $shape.Width = $Width
    $shape.Height = $Height
}

# 1) Full citation / source link along the bottom of the slide.
Add-CitationBox "TextBox 1" `
    "[1]https://towardsdatascience.com/taking-the-confusion-out-of-confusion-matrices-c1ce054b3d3e" `
    231.98725129448817 518.1890564181102 808.2856750913386 21.810944881889764

# 2) "[1]" marker next to the confusion-matrix picture.
Add-CitationBox "TextBox 2" "[1]" `
    671.9999212598425 478.6890551181102 37.71425196850394 21.810944881889764

# 3) "[1]" marker next to the other figure higher up on the slide.
Add-CitationBox "TextBox 3" "[1]" `
    724.0129921259843 264.4961417322835 37.71425196850394 21.810944881889764
